# Generate Report for Handback
# Updates the handback-status report: the "50c07ff8..." file's UUID is
# replaced by "9e9f19c9-ffff-4679-adca-478f414c2997" and the
# "d9fc61df..." file's UUID is replaced by
# "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99", together with refreshed
# handoff/handback timestamps and xliff file names, across the
# Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.md"
$ws1.Range("B2").Value = "e2e\9e9f19c9-ffff-4679-adca-478f414c2997.md"
$ws1.Range("G2").Value = "2016-08-31 15:21:26"

$ws1.Range("A3").Value = "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md"
$ws1.Range("B3").Value = "e2e\ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md"
$ws1.Range("G3").Value = "2016-08-31 15:21:26"

# Hyperlinks: this runtime's Hyperlink objects can't be edited or
# deleted individually (Delete()/property writes on a single Item
# silently create a stray extra hyperlink) so rebuild the whole
# collection for the sheet, preserving the original target URLs.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9faa7ea3be01b92f08f1e675beb137f4c5a5fbc/e2e/50c07ff8-5a04-4730-b635-71d1044566d3.md", "", "", "e2e\9e9f19c9-ffff-4679-adca-478f414c2997.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9faa7ea3be01b92f08f1e675beb137f4c5a5fbc/e2e/d9fc61df-b1b2-4fc7-908f-0debc97204aa.md", "", "", "e2e\ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.md"
$ws2.Range("G2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-31 15:21:22"
$ws2.Range("I2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.md"
$ws2.Range("J2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-31 15:21:48"

$ws2.Range("A3").Value = "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md"
$ws2.Range("G3").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-31 15:21:22"
$ws2.Range("I3").Value = "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md"
$ws2.Range("J3").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-31 15:21:48"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9faa7ea3be01b92f08f1e675beb137f4c5a5fbc/e2e/50c07ff8-5a04-4730-b635-71d1044566d3.md", "", "", "9e9f19c9-ffff-4679-adca-478f414c2997.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2d89b6ba67b74fbc1d0e3094d5981fbd2c726a7f/e2e/50c07ff8-5a04-4730-b635-71d1044566d3.md", "", "", "9e9f19c9-ffff-4679-adca-478f414c2997.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9faa7ea3be01b92f08f1e675beb137f4c5a5fbc/e2e/d9fc61df-b1b2-4fc7-908f-0debc97204aa.md", "", "", "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2d89b6ba67b74fbc1d0e3094d5981fbd2c726a7f/e2e/d9fc61df-b1b2-4fc7-908f-0debc97204aa.md", "", "", "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md")

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.md"
$ws3.Range("G2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-31 15:21:26"
$ws3.Range("I2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.md"
$ws3.Range("J2").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-31 15:21:56"

$ws3.Range("A3").Value = "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md"
$ws3.Range("G3").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-31 15:21:26"
$ws3.Range("I3").Value = "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md"
$ws3.Range("J3").Value = "9e9f19c9-ffff-4679-adca-478f414c2997.71d6acbb8c881ce9720bd11edcddd0a0d621658d.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-31 15:21:56"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9faa7ea3be01b92f08f1e675beb137f4c5a5fbc/e2e/50c07ff8-5a04-4730-b635-71d1044566d3.md", "", "", "9e9f19c9-ffff-4679-adca-478f414c2997.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/13fae09f2b0cda11ed9c250723d1b0db6abbfb99/e2e/50c07ff8-5a04-4730-b635-71d1044566d3.md", "", "", "9e9f19c9-ffff-4679-adca-478f414c2997.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a9faa7ea3be01b92f08f1e675beb137f4c5a5fbc/e2e/d9fc61df-b1b2-4fc7-908f-0debc97204aa.md", "", "", "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/13fae09f2b0cda11ed9c250723d1b0db6abbfb99/e2e/d9fc61df-b1b2-4fc7-908f-0debc97204aa.md", "", "", "ffff18cb22b3-27ed-4e28-be03-bac02a21ff99.md")
